# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# worksheet with the latest scrape values, row by row (rows 2..51).
#
# The source feed renders both columns as plain text (prices keep the
# site's "thousands-dot" grouping like "28.048.45", and the volume column
# keeps its padding spaces around the percentage, e.g. "  -3.94%  "), so
# every write below is forced to stay text even when a price looks like a
# plain decimal number (e.g. "226.66") that Excel would otherwise swallow
# into a numeric cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = '28.048.45'; E = '  -3.94%  ' },
    @{ Row = 3; D = '1.744.23'; E = '  -4.53%  ' },
    @{ Row = 4; D = $null; E = '  -0.15%  ' },
    @{ Row = 5; D = '226.66'; E = '  -3.43%  ' },
    @{ Row = 6; D = '0.5793'; E = '  -3.59%  ' },
    @{ Row = 7; D = $null; E = '  -0.11%  ' },
    @{ Row = 8; D = '0.2714'; E = '  -1.38%  ' },
    @{ Row = 9; D = '23.17'; E = '  -1.08%  ' },
    @{ Row = 10; D = '0.06595'; E = '  -5.05%  ' },
    @{ Row = 11; D = '0.07525'; E = '  -0.86%  ' },
    @{ Row = 12; D = '1.743.00'; E = '  -4.89%  ' },
    @{ Row = 13; D = '4.723'; E = '  -0.23%  ' },
    @{ Row = 14; D = '0.6053'; E = '  -3.19%  ' },
    @{ Row = 15; D = '1.982.90'; E = '  -4.47%  ' },
    @{ Row = 16; D = '74.33'; E = '  -3.95%  ' },
    @{ Row = 17; D = '0.000008690'; E = '  -11.59%  ' },
    @{ Row = 18; D = '28.041.99'; E = '  -2.91%  ' },
    @{ Row = 19; D = '5.332'; E = '  -4.33%  ' },
    @{ Row = 20; D = $null; E = '  -0.12%  ' },
    @{ Row = 21; D = '205.27'; E = '  -5.16%  ' },
    @{ Row = 22; D = $null; E = '  -2.20%  ' },
    @{ Row = 23; D = '6.636'; E = '  -3.69%  ' },
    @{ Row = 24; D = $null; E = '  -0.10%  ' },
    @{ Row = 25; D = '149.98'; E = '  -3.94%  ' },
    @{ Row = 26; D = '8.067'; E = '  +1.65%  ' },
    @{ Row = 27; D = $null; E = '  -4.24%  ' },
    @{ Row = 28; D = $null; E = '  -2.14%  ' },
    @{ Row = 29; D = '0.06222'; E = '  -5.52%  ' },
    @{ Row = 30; D = '1.387'; E = '  -1.55%  ' },
    @{ Row = 31; D = '1.393'; E = '  -3.29%  ' },
    @{ Row = 32; D = '3.745'; E = '  -2.16%  ' },
    @{ Row = 33; D = '3.713'; E = '  -1.60%  ' },
    @{ Row = 34; D = '1.678'; E = '  -2.82%  ' },
    @{ Row = 35; D = '1.037'; E = '  -5.12%  ' },
    @{ Row = 36; D = '0.6383'; E = '  -1.35%  ' },
    @{ Row = 37; D = '2.450'; E = '  -3.31%  ' },
    @{ Row = 38; D = '2.730'; E = '  -0.62%  ' },
    @{ Row = 39; D = '0.01673'; E = '  -4.93%  ' },
    @{ Row = 40; D = '1.128.90'; E = '  -1.33%  ' },
    @{ Row = 41; D = '6.201'; E = '  -4.35%  ' },
    @{ Row = 42; D = '0.8737'; E = '  -1.58%  ' },
    @{ Row = 43; D = '1.005'; E = '  +0.16%  ' },
    @{ Row = 44; D = '99.64'; E = '  -0.63%  ' },
    @{ Row = 45; D = '1.894.48'; E = '  -4.77%  ' },
    @{ Row = 46; D = '59.47'; E = '  -3.76%  ' },
    @{ Row = 47; D = '1.582'; E = $null },
    @{ Row = 48; D = $null; E = '  -5.57%  ' },
    @{ Row = 49; D = '8.295'; E = '  -2.32%  ' },
    @{ Row = 50; D = '0.05377'; E = '  -2.29%  ' },
    @{ Row = 51; D = '6.274'; E = '  -2.09%  ' }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($u.Row, 4)   # column D = Price
        if ($u.D -match '^[0-9]+(\.[0-9]+)?$') {
            # Looks like a plain decimal (e.g. "226.66") - Excel would
            # otherwise coerce the assignment into a Number cell, so pin
            # the format to Text first to keep the original string intact.
            $cell.NumberFormat = "@"
        }
        $cell.Value = $u.D
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E   # column E = Volume(1h)
    }
}
